$wb = $excel.ActiveWorkbook

# Update the CO2 price values on each year's sheet.
$wb.Worksheets.Item("2025").Range("A2").Value = 210
$wb.Worksheets.Item("2030").Range("A2").Value = 230
$wb.Worksheets.Item("2035").Range("A2").Value = 250
$wb.Worksheets.Item("2040").Range("A2").Value = 270
$wb.Worksheets.Item("2045").Range("A2").Value = 290
$wb.Worksheets.Item("2050").Range("A2").Value = 308

# Make the "2025" sheet the active/selected tab (was "2050").
$wb.Worksheets.Item("2025").Activate()
